# Add two new columns (I: "I0", J: "IF") with per-row numeric values to sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new columns, matching the existing header row styling
# (copy the "IP" header's formatting onto the two new header cells, then set text).
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)
$ws.Cells.Item(1, 10).PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Row data: (row, I0 value, IF value)
$data = @(
    @(2,9,9),
    @(3,8,8),
    @(4,8,8),
    @(5,7,7),
    @(6,7,7),
    @(7,8,8),
    @(8,9,9),
    @(9,9,9),
    @(10,7,7),
    @(11,9,9),
    @(12,9,9),
    @(13,7,7),
    @(14,8,8),
    @(15,6,6),
    @(16,9,9),
    @(17,5,5),
    @(18,7,7),
    @(19,7,8),
    @(20,8,8),
    @(21,8,8),
    @(22,8,8),
    @(23,8,8),
    @(24,8,8),
    @(25,9,9),
    @(26,7,7),
    @(27,7,7),
    @(28,3,4),
    @(29,8,8),
    @(30,6,6),
    @(31,9,9),
    @(32,9,9),
    @(33,9,9),
    @(34,8,8),
    @(35,8,8),
    @(36,7,7),
    @(37,6,6),
    @(38,7,7),
    @(39,7,8),
    @(40,6,6),
    @(41,8,8),
    @(42,8,8),
    @(43,9,9),
    @(44,8,8),
    @(45,6,6),
    @(46,10,10),
    @(47,8,8),
    @(48,9,9),
    @(49,8,8),
    @(50,7,7),
    @(51,8,8),
    @(52,6,7),
    @(53,8,8),
    @(54,8,8),
    @(55,7,7),
    @(56,8,8),
    @(57,6,7),
    @(58,6,7),
    @(59,7,7),
    @(60,7,7),
    @(61,5,5),
    @(62,7,7),
    @(63,7,7),
    @(64,8,8),
    @(65,8,8),
    @(66,9,9),
    @(67,7,7),
    @(68,5,5),
    @(69,5,5)
)

foreach ($row in $data) {
    $r = $row[0]
    $i0 = $row[1]
    $iF = $row[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $iF
}
